# The slide titled "Requisite mathematician joke" is removed from the
# deck (author replaced the joke slide with real content on the
# following slides). Locate it defensively by title text, falling back
# to the known position (slide 10) if the title can't be matched.

$p = $ppt.ActivePresentation

$targetIndex = -1
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $s = $p.Slides.Item($i)
    foreach ($shape in $s.Shapes) {
        if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
            if ($shape.TextFrame.TextRange.Text -eq "Requisite mathematician joke") {
                $targetIndex = $i
            }
            break
        }
    }
    if ($targetIndex -ne -1) {
        break
    }
}

if ($targetIndex -eq -1) {
    $targetIndex = 10
}

$p.Slides.Item($targetIndex).Delete()
